{"js": "// Apply the \"Helios\" report edits:\n//  - Merge runs that were split around proofing-error marks (Helios, dashboard,\n//    Grafana, alerting) back into single plain runs (purely cosmetic cleanup,\n//    text itself is unchanged).\n//  - Update the project end date (31/06/2025 -> 30/06/2023).\n//  - Update the headcount (54 -> 23 personnes).\n//  - Update the budget figure (250k\u20ac -> 2 300 000 \u20ac).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Helper: replace the full text of the paragraph at `index` with `newText`,\n// collapsing any split runs / proofErr markers into a single clean run.\n// Clearing the paragraph's whole range first (rather than just inserting\n// with Replace) also sweeps away any trailing proofing-error markers that\n// sit at the very end of the paragraph (after the last run).\nfunction replaceParagraphText(index, newText) {\n  const range = paragraphs.items[index].getRange(\"Whole\");\n  range.clear();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nconst items = paragraphs.items;\n\n// 0: Title \u2014 \"Projet \" + \"Helios\" + \" \u2014 Synth\u00e8se Mensuelle (Ao\u00fbt 2025)\"\nreplaceParagraphText(0, \"Projet Helios \u2014 Synth\u00e8se Mensuelle (Ao\u00fbt 2025)\");\n\n// 3: Context paragraph \u2014 \"L\u2019\u00e9quipe \" + \"Helios\" + \" a \u00e9t\u00e9 mobilis\u00e9e ... juillet.\"\nreplaceParagraphText(\n  3,\n  \"L\u2019\u00e9quipe Helios a \u00e9t\u00e9 mobilis\u00e9e pour soutenir la transformation num\u00e9rique du client dans le secteur de la logistique. Leur besoin portait sur l\u2019am\u00e9lioration du suivi des exp\u00e9ditions en temps r\u00e9el, avec une forte contrainte d\u2019int\u00e9gration dans les syst\u00e8mes existants. La collaboration s\u2019est amorc\u00e9e fin mai, avec une mont\u00e9e en charge progressive jusqu\u2019en juillet.\"\n);\n\n// 9: \"- D\u00e9ploiement d\u2019un prototype de dashboard temps r\u00e9el bas\u00e9 sur Grafana.\"\nreplaceParagraphText(\n  9,\n  \"- D\u00e9ploiement d\u2019un prototype de dashboard temps r\u00e9el bas\u00e9 sur Grafana.\"\n);\n\n// 11: \"- Mise en place d\u2019un syst\u00e8me d\u2019alerting via Slack et SMS pour les anomalies critiques.\"\nreplaceParagraphText(\n  11,\n  \"- Mise en place d\u2019un syst\u00e8me d\u2019alerting via Slack et SMS pour les anomalies critiques.\"\n);\n\n// 18: \"Collaboration fluide ... notamment sur l\u2019ergonomie du dashboard.\"\nreplaceParagraphText(\n  18,\n  \"Collaboration fluide avec les \u00e9quipes IT du client. Le choix de technologies open-source a permis une adoption rapide et des \u00e9conomies substantielles. Les retours utilisateurs sont tr\u00e8s positifs, notamment sur l\u2019ergonomie du dashboard.\"\n);\n\n// 21: \"- Grafana\"\nreplaceParagraphText(21, \"- Grafana\");\n\n// 28: \"Le projet s\u2019est \u00e9tal\u00e9 du 02/04/2021 au 31/06/2025 et a impliqu\u00e9 54 personnes.\"\nreplaceParagraphText(\n  28,\n  \"Le projet s\u2019est \u00e9tal\u00e9 du 02/04/2021 au 30/06/2023 et a impliqu\u00e9 23 personnes.\"\n);\n\n// 30: \"Le budget total a \u00e9t\u00e9 de 250k\u20ac.\"\nreplaceParagraphText(30, \"Le budget total a \u00e9t\u00e9 de 2 300 000 \u20ac.\");\n\nawait context.sync();\n", "ps1": "# Apply the \"Helios\" report edits:\n#  - Merge runs that were split around proofing-error marks (Helios, dashboard,\n#    Grafana, alerting) back into single plain runs (purely cosmetic cleanup,\n#    the visible text itself is unchanged).\n#  - Update the project end date (31/06/2025 -> 30/06/2023).\n#  - Update the headcount (54 -> 23 personnes).\n#  - Update the budget figure (250k\u20ac -> 2 300 000 \u20ac).\n\n# Replace the full text of paragraph number $index (1-based, Word COM style)\n# with $newText. Deleting the paragraph's content together with its\n# paragraph mark (up to, but not including, the start of the following\n# paragraph) and then re-inserting fresh text + a paragraph mark guarantees\n# any left-over <w:proofErr/> markers (which Word leaves behind around\n# spell-checked words such as \"Helios\", \"dashboard\", \"Grafana\", \"alerting\")\n# are swept away as well, collapsing the paragraph down to a single clean\n# run \u2014 exactly like Word does when you retype the whole line.\nfunction Replace-ParagraphText {\n    param(\n        $doc,\n        [int]$index,\n        [string]$newText\n    )\n    $p = $doc.Paragraphs.Item($index)\n    $pNext = $doc.Paragraphs.Item($index + 1)\n    $r = $doc.Range($p.Range.Start, $pNext.Range.Start)\n    $r.Delete()\n    $ins = $doc.Range($p.Range.Start, $p.Range.Start)\n    $ins.InsertBefore($newText + \"`r\")\n}\n\n$d = $word.ActiveDocument\n\n# Paragraph 1: title \u2014 \"Projet \" + \"Helios\" + \" \u2014 Synth\u00e8se Mensuelle (Ao\u00fbt 2025)\"\nReplace-ParagraphText $d 1 \"Projet Helios \u2014 Synth\u00e8se Mensuelle (Ao\u00fbt 2025)\"\n\n# Paragraph 4: \"L'\u00e9quipe \" + \"Helios\" + \" a \u00e9t\u00e9 mobilis\u00e9e ... juillet.\"\nReplace-ParagraphText $d 4 \"L\u2019\u00e9quipe Helios a \u00e9t\u00e9 mobilis\u00e9e pour soutenir la transformation num\u00e9rique du client dans le secteur de la logistique. Leur besoin portait sur l\u2019am\u00e9lioration du suivi des exp\u00e9ditions en temps r\u00e9el, avec une forte contrainte d\u2019int\u00e9gration dans les syst\u00e8mes existants. La collaboration s\u2019est amorc\u00e9e fin mai, avec une mont\u00e9e en charge progressive jusqu\u2019en juillet.\"\n\n# Paragraph 10: \"- D\u00e9ploiement d'un prototype de dashboard temps r\u00e9el bas\u00e9 sur Grafana.\"\nReplace-ParagraphText $d 10 \"- D\u00e9ploiement d\u2019un prototype de dashboard temps r\u00e9el bas\u00e9 sur Grafana.\"\n\n# Paragraph 12: \"- Mise en place d'un syst\u00e8me d'alerting via Slack et SMS pour les anomalies critiques.\"\nReplace-ParagraphText $d 12 \"- Mise en place d\u2019un syst\u00e8me d\u2019alerting via Slack et SMS pour les anomalies critiques.\"\n\n# Paragraph 19: \"Collaboration fluide ... notamment sur l'ergonomie du dashboard.\"\nReplace-ParagraphText $d 19 \"Collaboration fluide avec les \u00e9quipes IT du client. Le choix de technologies open-source a permis une adoption rapide et des \u00e9conomies substantielles. Les retours utilisateurs sont tr\u00e8s positifs, notamment sur l\u2019ergonomie du dashboard.\"\n\n# Paragraph 22: \"- Grafana\"\nReplace-ParagraphText $d 22 \"- Grafana\"\n\n# Paragraph 29: project dates / headcount \u2014 31/06/2025 -> 30/06/2023, 54 -> 23 personnes\nReplace-ParagraphText $d 29 \"Le projet s\u2019est \u00e9tal\u00e9 du 02/04/2021 au 30/06/2023 et a impliqu\u00e9 23 personnes.\"\n\n# Paragraph 31: budget \u2014 250k\u20ac -> 2 300 000 \u20ac\nReplace-ParagraphText $d 31 \"Le budget total a \u00e9t\u00e9 de 2 300 000 \u20ac.\"\n"}
